# Daily attendance processing - 2025-12-17 13:43:54
# Normalize the "Recorded By" (column G) attribution strings: reorder the
# comma-separated author list so the real user's email is listed before
# the automated "System"/"system" entries (and before "admin@admin.com").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$col = $ws.Range("G1:G157")

$col.Replace("backup@backdoor.com, System, system", "backup@backdoor.com, system, System")
$col.Replace("System, dnasr281@gmail.com", "dnasr281@gmail.com, System")
$col.Replace("admin@admin.com, dnasr281@gmail.com", "dnasr281@gmail.com, admin@admin.com")
